# Auto-generated edit script for Chocobo_Profits workbook
# Updates currentAveragePrice / Leve price / profit columns (H,I,J,K,L,M,N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 11799.9
$ws.Range("I21").Value = 7600
$ws.Range("K21").Value = 7600
$ws.Range("M21").Value = -7132
# Row 23
$ws.Range("H23").Value = 11799.9
$ws.Range("I23").Value = 7600
$ws.Range("K23").Value = 7600
$ws.Range("M23").Value = -7366
# Row 33
$ws.Range("H33").Value = 184.16667
$ws.Range("I33").Value = 164.45454
$ws.Range("K33").Value = 164.45454
$ws.Range("M33").Value = 64.54545999999999
# Row 53
$ws.Range("H53").Value = 525.7
$ws.Range("I53").Value = 278
$ws.Range("J53").Value = 985.7143
$ws.Range("K53").Value = 278
$ws.Range("L53").Value = 985.7143
$ws.Range("M53").Value = 359
$ws.Range("N53").Value = -2259.7143
# Row 98
$ws.Range("H98").Value = 5020.9463
$ws.Range("I98").Value = 3101.9
$ws.Range("J98").Value = 7235.231
$ws.Range("K98").Value = 3101.9
$ws.Range("L98").Value = 7235.231
$ws.Range("M98").Value = -1603.9
$ws.Range("N98").Value = -10231.231
# Row 108
$ws.Range("H108").Value = 35971.43
$ws.Range("J108").Value = 35971.43
$ws.Range("L108").Value = 35971.43
$ws.Range("N108").Value = -43651.43
# Row 120
$ws.Range("H120").Value = 28995
$ws.Range("J120").Value = 28995
$ws.Range("L120").Value = 28995
$ws.Range("N120").Value = -38671
# Row 122
$ws.Range("H122").Value = 5020.9463
$ws.Range("I122").Value = 3101.9
$ws.Range("J122").Value = 7235.231
$ws.Range("K122").Value = 9305.700000000001
$ws.Range("L122").Value = 21705.693
$ws.Range("M122").Value = -6855.700000000001
$ws.Range("N122").Value = -26605.693
# Row 129
$ws.Range("H129").Value = 836.96
$ws.Range("I129").Value = 345.5
$ws.Range("J129").Value = 857.4375
$ws.Range("K129").Value = 1036.5
$ws.Range("L129").Value = 2572.3125
$ws.Range("M129").Value = 3963.5
$ws.Range("N129").Value = -12572.3125
# Row 133
$ws.Range("H133").Value = 39608.184
$ws.Range("J133").Value = 39608.184
$ws.Range("L133").Value = 39608.184
$ws.Range("N133").Value = -49728.184
# Row 137
$ws.Range("H137").Value = 4765519
$ws.Range("I137").Value = 6803627
$ws.Range("J137").Value = 9933.333000000001
$ws.Range("K137").Value = 20410881
$ws.Range("L137").Value = 29799.999
$ws.Range("M137").Value = -20408331
$ws.Range("N137").Value = -34899.999
# Row 139
$ws.Range("H139").Value = 45333.57
$ws.Range("J139").Value = 45333.57
$ws.Range("L139").Value = 45333.57
$ws.Range("N139").Value = -55613.57

# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 32
$ws.Range("H32").Value = 5714.5186
$ws.Range("I32").Value = 6779.483
$ws.Range("J32").Value = 4479.16
$ws.Range("K32").Value = 6779.483
$ws.Range("L32").Value = 4479.16
$ws.Range("M32").Value = -6492.483
$ws.Range("N32").Value = -5053.16
# Row 74
$ws.Range("H74").Value = 2545.25
$ws.Range("I74").Value = 1566.1666
$ws.Range("J74").Value = 5482.5
$ws.Range("K74").Value = 1566.1666
$ws.Range("L74").Value = 5482.5
$ws.Range("M74").Value = -692.1666
$ws.Range("N74").Value = -7230.5
# Row 77
$ws.Range("H77").Value = 2545.25
$ws.Range("I77").Value = 1566.1666
$ws.Range("J77").Value = 5482.5
$ws.Range("K77").Value = 7830.833000000001
$ws.Range("L77").Value = 27412.5
$ws.Range("M77").Value = -3462.833000000001
$ws.Range("N77").Value = -36148.5
# Row 115
$ws.Range("H115").Value = 26888
$ws.Range("J115").Value = 26888
$ws.Range("L115").Value = 26888
$ws.Range("N115").Value = -30022
# Row 117
$ws.Range("H117").Value = 27764
$ws.Range("J117").Value = 27764
$ws.Range("L117").Value = 27764
$ws.Range("N117").Value = -36942
# Row 119
$ws.Range("H119").Value = 34638.332
$ws.Range("J119").Value = 34638.332
$ws.Range("L119").Value = 34638.332
$ws.Range("N119").Value = -44314.332
# Row 122
$ws.Range("H122").Value = 4151.8335
$ws.Range("I122").Value = 3597.52
$ws.Range("J122").Value = 5411.636
$ws.Range("K122").Value = 10792.56
$ws.Range("L122").Value = 16234.908
$ws.Range("M122").Value = -8342.559999999999
$ws.Range("N122").Value = -21134.908
# Row 139
$ws.Range("H139").Value = 42132.812
$ws.Range("J139").Value = 42132.812
$ws.Range("L139").Value = 42132.812
$ws.Range("N139").Value = -52412.812

# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")
# Row 9
$ws.Range("H9").Value = 13332.667
$ws.Range("J9").Value = 13332.667
$ws.Range("L9").Value = 13332.667
$ws.Range("N9").Value = -13668.667

# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5278.476
$ws.Range("I31").Value = 2084.4324
$ws.Range("J31").Value = 9823.846
$ws.Range("K31").Value = 2084.4324
$ws.Range("L31").Value = 9823.846
$ws.Range("M31").Value = -1789.4324
$ws.Range("N31").Value = -10413.846
# Row 34
$ws.Range("H34").Value = 5278.476
$ws.Range("I34").Value = 2084.4324
$ws.Range("J34").Value = 9823.846
$ws.Range("K34").Value = 2084.4324
$ws.Range("L34").Value = 9823.846
$ws.Range("M34").Value = -1882.4324
$ws.Range("N34").Value = -10227.846
# Row 59
$ws.Range("H59").Value = 33797.875
$ws.Range("J59").Value = 33797.875
$ws.Range("L59").Value = 33797.875
$ws.Range("N59").Value = -36087.875
# Row 99
$ws.Range("H99").Value = 10530361
$ws.Range("I99").Value = 16668239
$ws.Range("J99").Value = 8284.286
$ws.Range("K99").Value = 16668239
$ws.Range("L99").Value = 8284.286
$ws.Range("M99").Value = -16666741
$ws.Range("N99").Value = -11280.286
# Row 107
$ws.Range("H107").Value = 775.65216
$ws.Range("I107").Value = 607.7778
$ws.Range("J107").Value = 1380
$ws.Range("K107").Value = 607.7778
$ws.Range("L107").Value = 1380
$ws.Range("M107").Value = 1312.2222
$ws.Range("N107").Value = -5220
# Row 126
$ws.Range("H126").Value = 10530361
$ws.Range("I126").Value = 16668239
$ws.Range("J126").Value = 8284.286
$ws.Range("K126").Value = 50004717
$ws.Range("L126").Value = 24852.858
$ws.Range("M126").Value = -50002247
$ws.Range("N126").Value = -29792.858

# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 743296.5
$ws.Range("I5").Value = 572
$ws.Range("K5").Value = 1716
$ws.Range("M5").Value = -1604
# Row 122
$ws.Range("H122").Value = 2659.537
$ws.Range("J122").Value = 2930.3696
$ws.Range("L122").Value = 26373.3264
$ws.Range("N122").Value = -31273.3264
# Row 131
$ws.Range("H131").Value = 781.84
$ws.Range("J131").Value = 828.3955999999999
$ws.Range("L131").Value = 2485.1868
$ws.Range("N131").Value = -12565.1868
# Row 135
$ws.Range("H135").Value = 743296.5
$ws.Range("I135").Value = 572
$ws.Range("K135").Value = 5148
$ws.Range("M135").Value = -2613

# ==== Sheet: GSM ====
$ws = $wb.Worksheets.Item("GSM")
# Row 94
$ws.Range("H94").Value = 35000
$ws.Range("J94").Value = 35000
$ws.Range("L94").Value = 35000
$ws.Range("N94").Value = -36352
# Row 102
$ws.Range("H102").Value = 1913.8462
$ws.Range("I102").Value = 1352
$ws.Range("J102").Value = 2917.1428
$ws.Range("K102").Value = 1352
$ws.Range("L102").Value = 2917.1428
$ws.Range("M102").Value = 270
$ws.Range("N102").Value = -6161.1428
# Row 107
$ws.Range("H107").Value = 7408036
$ws.Range("I107").Value = 380.55554
$ws.Range("J107").Value = 18519520
$ws.Range("K107").Value = 380.55554
$ws.Range("L107").Value = 18519520
$ws.Range("M107").Value = 1539.44446
$ws.Range("N107").Value = -18523360
# Row 122
$ws.Range("H122").Value = 4352.8667
$ws.Range("I122").Value = 4089.125
$ws.Range("J122").Value = 4654.2856
$ws.Range("K122").Value = 12267.375
$ws.Range("L122").Value = 13962.8568
$ws.Range("M122").Value = -9817.375
$ws.Range("N122").Value = -18862.8568
# Row 132
$ws.Range("H132").Value = 5999.5713
$ws.Range("I132").Value = 1999.5
$ws.Range("K132").Value = 5998.5
$ws.Range("M132").Value = -3468.5
# Row 138
$ws.Range("H138").Value = 41772.25
$ws.Range("J138").Value = 41772.25
$ws.Range("L138").Value = 41772.25
$ws.Range("N138").Value = -52052.25

# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 7123.375
$ws.Range("I2").Value = 989.5
$ws.Range("J2").Value = 9168
$ws.Range("K2").Value = 989.5
$ws.Range("L2").Value = 9168
$ws.Range("M2").Value = -877.5
$ws.Range("N2").Value = -9392
# Row 55
$ws.Range("H55").Value = 85.70587999999999
$ws.Range("I55").Value = 70
$ws.Range("J55").Value = 136.75
$ws.Range("K55").Value = 70
$ws.Range("L55").Value = 136.75
$ws.Range("M55").Value = 103
$ws.Range("N55").Value = -482.75
# Row 132
$ws.Range("H132").Value = 4012.5
$ws.Range("I132").Value = 3065.75
$ws.Range("K132").Value = 9197.25
$ws.Range("M132").Value = -6667.25
# Row 141
$ws.Range("H141").Value = 31918.572
$ws.Range("J141").Value = 31918.572
$ws.Range("L141").Value = 31918.572
$ws.Range("N141").Value = -42278.572

# ==== Sheet: WVR ====
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 4169.273
$ws.Range("I122").Value = 3210.9375
$ws.Range("K122").Value = 9632.8125
$ws.Range("M122").Value = -7182.8125
# Row 141
$ws.Range("H141").Value = 43686
$ws.Range("J141").Value = 43686
$ws.Range("L141").Value = 43686
$ws.Range("N141").Value = -54046

